$d = $word.ActiveDocument

$target = "Il est vrai que nos universit" + [char]0x00E9 + "s n" + [char]0x2019 + "ayant pas dans leur cursus les derni" + [char]0x00E8 + "res mises " + [char]0x00E0 + " jour concernant l" + [char]0x2019 + [char]0x00E9 + "volution des m" + [char]0x00E9 + "tiers de d" + [char]0x00E9 + "veloppeurs, ne nous pr" + [char]0x00E9 + "parent absolument pas du tout pour le march" + [char]0x00E9 + " de l" + [char]0x2019 + "emploi, et je comprends facilement que des " + [char]0x00E9 + "tudiants en finissant se sentent un peu d" + [char]0x00E9 + "rout" + [char]0x00E9 + "s. "

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text.TrimEnd() -eq $target.TrimEnd()) {
        $r.Font.Size = 11
        break
    }
}
